# Update countries & provincias Spain
# Applies the daily data refresh + re-sort of the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 00:35"

# --- 2. Update totals for Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 1723638   # Casos totales
$ws.Cells.Item(4, 3).Value = 17412     # Nuevos casos
$ws.Cells.Item(4, 4).Value = 477314    # Casos activos
$ws.Cells.Item(4, 5).Value = 1145827   # Recuperados
$ws.Cells.Item(4, 7).Value = 692       # Casos criticos
$ws.Cells.Item(4, 8).Value = 100497    # Muertes

# --- 3. Update totals for Camerun (row 68) ---
$ws.Cells.Item(68, 2).Value = 5436
$ws.Cells.Item(68, 3).Value = 546
$ws.Cells.Item(68, 4).Value = 1996
$ws.Cells.Item(68, 5).Value = 3265
$ws.Cells.Item(68, 7).Value = 10
$ws.Cells.Item(68, 8).Value = 175

# --- 4. Update totals for Guinea Ecuatorial (row 144) ---
$ws.Cells.Item(144, 2).Value = 384
$ws.Cells.Item(144, 3).Value = 31
$ws.Cells.Item(144, 5).Value = 237

# --- 5. Bermudas overtakes Islas Caimanes: swap the two rows (name + stats) ---
# Row 165 held Islas Caimanes, row 166 held Bermudas (sorted by Casos totales desc).
# Bermudas' count grew above Islas Caimanes', so the rows swap places while
# Bermudas also receives its updated daily figures.
$ws.Cells.Item(165, 1).Value = "Bermudas"
$ws.Cells.Item(165, 2).Value = 139
$ws.Cells.Item(165, 3).Value = 6
$ws.Cells.Item(165, 4).Value = 89
$ws.Cells.Item(165, 5).Value = 41
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 9

$ws.Cells.Item(166, 1).Value = "Islas Caimanes"
$ws.Cells.Item(166, 2).Value = 134
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 61
$ws.Cells.Item(166, 5).Value = 72
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 1

# --- 6. Islas Turcas y Caicos overtakes Groenlandia (tie-break reorder) ---
# Row 207 held Groenlandia, row 208 held Islas Turcas y Caicos; both have the
# same Casos totales (12), but the tie-break order changed, so the rows swap.
$ws.Cells.Item(207, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(207, 2).Value = 12
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 10
$ws.Cells.Item(207, 5).Value = 1
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 1

$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 5).Value = 1
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0
